$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For every row on the sheet, if the Venue (column C) is "P&G",
# shift the Time (column F) back by 30 minutes (0.5/24 of a day).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $venue = $ws.Cells.Item($r, 3).Value2
    if ($venue -eq "P&G") {
        $timeCell = $ws.Cells.Item($r, 6)
        $totalMinutes = [Math]::Round($timeCell.Value2 * 1440)
        $timeCell.Value2 = ($totalMinutes - 30) / 1440
    }
}
